# Apply hybrid bold + color highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) in the achievements /
# work-experience bullet points, matching the commit's DOCX processing
# ("Added hybrid bold + color highlighting for impact metrics").
#
# Each target run of text is split into: plain-text run(s) + a bold,
# colored (#2C3E50) run for each numeric metric, by scoping Find.Execute
# to the owning paragraph's Range and toggling Font.Bold / Font.Color on
# the found sub-range for every metric token, left to right.

$d = $word.ActiveDocument

# #2C3E50 as an OLE_COLOR (BGR-packed) value: R=0x2C G=0x3E B=0x50
$metricColor = 5258796

function Highlight-Metrics {
    param(
        [int]$ParagraphIndex,
        [string[]]$Terms
    )

    $para = $d.Paragraphs.Item($ParagraphIndex)
    $paraStart = $para.Range.Start
    $paraEnd = $para.Range.End

    foreach ($term in $Terms) {
        $searchRange = $d.Range($paraStart, $paraEnd)
        $found = $searchRange.Find.Execute($term, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $searchRange.Font.Bold = $true
            $searchRange.Font.Color = $metricColor
        }
    }
}

# "• Discovered systematic race coding errors ... from 23% to 64%"
Highlight-Metrics 9 @('23%', '64%')

# "• Achieved 87% prediction accuracy ... standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"
Highlight-Metrics 11 @('87%', '71%', [char]0x00B1 + '4.2%', [char]0x00B1 + '2.1%')

# "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
Highlight-Metrics 31 @('1,200')

# "• Created comprehensive meta-analysis framework ... became the $400M Polling Consortium Database ... now valued at $1B+"
Highlight-Metrics 46 @('$400M', '$1B')

# "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Highlight-Metrics 63 @('73.5%', '$4.7M')

# "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
Highlight-Metrics 65 @('87%', '71%')

Write-Output "Applied quantitative metric highlighting to 6 paragraphs."
